$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($r1, $r2, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value2 = $v2
        $ws.Cells.Item($r2, $c).Value2 = $v1
    }
}

# Swap match data (columns B:AC) between row pairs whose sort order changed.
# Column A (the running id counter) stays put on each row.
Swap-RowRange 100 101 2 29
Swap-RowRange 102 103 2 29
Swap-RowRange 162 163 2 29

# Row 183: updated odds for the existing fixture (Portuguesa vs Carabobo),
# which also moved up from the old row 186 slot.
$ws.Cells.Item(183, 2).Value2  = 7977863
$ws.Cells.Item(183, 5).Value2  = 45381.89583333334
$ws.Cells.Item(183, 6).Value2  = "Portuguesa"
$ws.Cells.Item(183, 7).Value2  = "Carabobo"
$ws.Cells.Item(183, 11).Value2 = 2.75
$ws.Cells.Item(183, 12).Value2 = 2.875
$ws.Cells.Item(183, 13).Value2 = 2.55
$ws.Cells.Item(183, 14).Value2 = 3.2
$ws.Cells.Item(183, 15).Value2 = 2.75
$ws.Cells.Item(183, 16).Value2 = 2.375
$ws.Cells.Item(183, 17).Value2 = 0.25
$ws.Cells.Item(183, 18).Value2 = 1.775
$ws.Cells.Item(183, 19).Value2 = 2.025
$ws.Cells.Item(183, 20).Value2 = 2

# Row 184: new fixture (Estudiantes Merida vs Deportivo La Guaira) replacing the old one.
$ws.Cells.Item(184, 2).Value2  = 7977380
$ws.Cells.Item(184, 5).Value2  = 45382.70833333334
$ws.Cells.Item(184, 6).Value2  = "Estudiantes Merida"
$ws.Cells.Item(184, 7).Value2  = "Deportivo La Guaira"
$ws.Cells.Item(184, 11).Value2 = 3.4
$ws.Cells.Item(184, 13).Value2 = 2.15
$ws.Cells.Item(184, 14).Value2 = 3.25
$ws.Cells.Item(184, 16).Value2 = 2.2
$ws.Cells.Item(184, 17).Value2 = 0.25
$ws.Cells.Item(184, 18).Value2 = 1.85
$ws.Cells.Item(184, 19).Value2 = 1.95
$ws.Cells.Item(184, 21).Value2 = 1.9
$ws.Cells.Item(184, 22).Value2 = 1.9

# Row 185: new fixture (Academia Puerto Cabello vs Monagas) replacing the old one.
$ws.Cells.Item(185, 2).Value2  = 8027687
$ws.Cells.Item(185, 5).Value2  = 45382.83333333334
$ws.Cells.Item(185, 6).Value2  = "Academia Puerto Cabello"
$ws.Cells.Item(185, 7).Value2  = "Monagas"
$ws.Cells.Item(185, 11).Value2 = 1.75
$ws.Cells.Item(185, 12).Value2 = 3.4
$ws.Cells.Item(185, 13).Value2 = 4.333
$ws.Cells.Item(185, 14).Value2 = 1.909
$ws.Cells.Item(185, 15).Value2 = 2.875
$ws.Cells.Item(185, 16).Value2 = 4.2
$ws.Cells.Item(185, 17).Value2 = -0.5
$ws.Cells.Item(185, 18).Value2 = 1.95
$ws.Cells.Item(185, 19).Value2 = 1.85
$ws.Cells.Item(185, 20).Value2 = 2
$ws.Cells.Item(185, 21).Value2 = 1.9
$ws.Cells.Item(185, 22).Value2 = 1.9

# Row 186 (old Portuguesa vs Carabobo row, now folded into row 183 above) is removed entirely.
$ws.Rows.Item(186).Delete()
